# This script rewrites the "system_info.docx" source-code block:
#   - the first paragraph (import statements) becomes style "FirstParagraph"
#     and loses its Pandoc syntax-highlighting run styles
#   - the print() calls for Operating System / Python Version / Machine /
#     Processor are split into their own "BodyText" paragraph
#   - the mem = ... / print(Total Memory) / print(Available Memory) lines
#     become a third "BodyText" paragraph, with the tail starting at
#     "2), 2)) print("Available Memory (MB):" ... (1024" made bold
#   - the verbatim command-output paragraph at the end is removed
#   - straight quotes in the printed strings become curly/smart quotes

$d = $word.ActiveDocument

$p4 = $d.Paragraphs(4)
$p5 = $d.Paragraphs(5)
$rng = $d.Range($p4.Range.Start, $p5.Range.End)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">import platform</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">import psutil</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">print(</w:t></w:r><w:r><w:t xml:space="preserve">&#x201c;Operating System:&#x201d;</w:t></w:r><w:r><w:t xml:space="preserve">, platform.system(), platform.release())</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">print(</w:t></w:r><w:r><w:t xml:space="preserve">&#x201c;Python Version:&#x201d;</w:t></w:r><w:r><w:t xml:space="preserve">, platform.python_version())</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">print(</w:t></w:r><w:r><w:t xml:space="preserve">&#x201c;Machine:&#x201d;</w:t></w:r><w:r><w:t xml:space="preserve">, platform.machine())</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">print(</w:t></w:r><w:r><w:t xml:space="preserve">&#x201c;Processor:&#x201d;</w:t></w:r><w:r><w:t xml:space="preserve">, platform.processor())</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr><w:r><w:t xml:space="preserve">mem = psutil.virtual_memory()</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">print(</w:t></w:r><w:r><w:t xml:space="preserve">&#x201c;Total Memory (MB):&#x201d;</w:t></w:r><w:r><w:t xml:space="preserve">, round(mem.total / (1024</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">2), 2))</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">print(</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">&#x201c;Available Memory (MB):&#x201d;</w:t></w:r><w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">, round(mem.available / (1024</w:t></w:r><w:r><w:t xml:space="preserve">2), 2))</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$rng.InsertXML($xml)

Write-Output "Paragraph count after edit:"
Write-Output $d.Paragraphs.Count
